# RubricaPlitix_SQLWarriors.xlsx - fill in evidence/status for the
# "seguridad" and "paquetes PL/SQL" sections of the rubric, flipping the
# related answers from NO to SI now that the work is done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Permisos section (rows 14-18): evidence column was empty, now points
# to the dba_roles view.
$ws.Range("E14").Value = "Consultar los roles creados por el usuario dba_roles"
$ws.Range("E15").Value = "Consultar los roles creados por el usuario dba_roles"
$ws.Range("E16").Value = "Consultar los roles creados por el usuario dba_roles"
$ws.Range("E17").Value = "Consultar los roles creados por el usuario dba_roles"
$ws.Range("E18").Value = "Consultar los roles creados por el usuario dba_roles"

# --- Transacciones / Excepciones (rows 30-31): add evidence.
$ws.Range("E30").Value = "Añadir un commit y un rollback en los procedimientos"
$ws.Range("E31").Value = "Mostrar en el código del paquete PL/SQL"

# --- Row 32: funcionalidades probadas -> ahora SI, con evidencia.
$ws.Range("D32").Value = "SI"
$ws.Range("E32").Value = "Ejecución de test con datos coherentes y verificar sus datos"

# --- Seguridad (rows 33-38): roles, usuarios, permisos, contraseñas, TDE.
$ws.Range("D34").Value = "SI"
$ws.Range("E34").Value = "Consultar dba_roles_privs"

$ws.Range("D35").Value = "SI"
$ws.Range("E35").Value = "Consultar dba_role_privs"

$ws.Range("D36").Value = "SI"
$ws.Range("E36").Value = "Consultar dba_role_privs"

$ws.Range("D37").Value = "SI"
$ws.Range("E37").Value = "Consultar dba_profiles"

$ws.Range("D38").Value = "SI"
$ws.Range("E38").Value = "Consultar dba_encrypted_columns"

# --- Modelo (rows 52-53): restricciones semánticas / NOT NULL-UNIQUE.
$ws.Range("D52").Value = "SI"
$ws.Range("E52").Value = "Comprobar que el modelo cumple con las restricciones semánticas (rangos válidos)"

$ws.Range("D53").Value = "SI"
$ws.Range("E53").Value = "Comprobar las columnas obligatorias y las claves únicas"

# --- Miscelanea (rows 55-57).
$ws.Range("D55").Value = "SI"
$ws.Range("E55").Value = "Verificar que existen productos asociados a cuentas válidas (COUNT > 0)"

$ws.Range("D56").Value = "SI"
$ws.Range("E56").Value = "Todos los objetos siguen la convención de nombres en minúsculas (COUNT = 0 para nombres con mayús.)"

$ws.Range("D57").Value = "SI"
$ws.Range("E57").Value = "Verificar los 3 objetos adicionales creados"

# --- Contexto (row 58): only the SI/NO flips, no evidence cell.
$ws.Range("D58").Value = "SI"
